$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29: Kazajistan -- refreshed totals ---
$ws.Range("B29").Value = 99442
$ws.Range("C29").Value = 741
$ws.Range("D29").Value = 72523
$ws.Range("E29").Value = 25861

# --- Rows 180-183: "Islas Turcas y Caicos" overtakes "Papua Nueva Guinea",
#     "Islas Caimanes" and "Gibraltar" in the case-count ranking, so the
#     sorted table reshuffles these four rows. Only Islas Turcas y Caicos's
#     figures actually change; the other three simply shift down one row
#     with their existing figures intact. ---
$ws.Range("A180").Value = "Islas Turcas y Caicos"
$ws.Range("B180").Value = 216
$ws.Range("C180").Value = 19
$ws.Range("D180").Value = 39
$ws.Range("E180").Value = 175
$ws.Range("F180").Value = 0
$ws.Range("G180").Value = 0
$ws.Range("H180").Value = 2

$ws.Range("A181").Value = "Papua Nueva Guinea"
$ws.Range("B181").Value = 214
$ws.Range("C181").Value = 0
$ws.Range("D181").Value = 53
$ws.Range("E181").Value = 158
$ws.Range("F181").Value = 0
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 3

$ws.Range("A182").Value = "Islas Caimanes"
$ws.Range("B182").Value = 203
$ws.Range("C182").Value = 0
$ws.Range("D182").Value = 202
$ws.Range("E182").Value = 0
$ws.Range("F182").Value = 0
$ws.Range("G182").Value = 0
$ws.Range("H182").Value = 1

$ws.Range("A183").Value = "Gibraltar"
$ws.Range("B183").Value = 201
$ws.Range("C183").Value = 0
$ws.Range("D183").Value = 186
$ws.Range("E183").Value = 15
$ws.Range("F183").Value = 0
$ws.Range("G183").Value = 0
$ws.Range("H183").Value = 0

# --- Row 191: Butan -- refreshed totals ---
$ws.Range("B191").Value = 110
$ws.Range("C191").Value = 2
$ws.Range("E191").Value = 14

# --- Rows 213-214: "Montserrat" and "Islas Malvinas" are tied on total
#     cases, and swap places in the sorted table (figures unchanged, just
#     the row each one occupies). ---
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

# --- Title cell: refresh "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 10 de Agosto de 2020 a las 06:09"
